$wb = $excel.ActiveWorkbook

$ov = $wb.Worksheets.Item("Overview")
$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

# Update the "Ready for handoff" status text to "In Translation" everywhere it appears.
$ov.Range("E2").Value = "In Translation"
$ov.Range("F2").Value = "In Translation"
$zh.Range("C2").Value = "In Translation"
$de.Range("C2").Value = "In Translation"

# Narrow the now-shorter "Status" columns to match the new content width.
$ov.Columns.Item(5).ColumnWidth = 12.5
$ov.Columns.Item(6).ColumnWidth = 12.5
$zh.Columns.Item(3).ColumnWidth = 12.5
$de.Columns.Item(3).ColumnWidth = 12.5
